$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.986.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.930.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "472.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.751"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000324"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.551.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.967.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.517.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("B25").Value = "EthereumClassic"
$ws.Range("C25").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "39.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.26%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "736.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.52%  "
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0485"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.39%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.69%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.341"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.142"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.66%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.57%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.91%  "
